$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E (shifts existing Arts/calle28/NCR/Delhi columns right)
$ws.Range("E1").EntireColumn.Insert()

# New header value: birth date text
$ws.Range("E1").Value = "27 Mar 1997"

# Update selection to match recorded state
$ws.Range("F9").Select() | Out-Null
